# Add duplicate detection for contract note imports:
# a new Buy transaction (2026-02-10) is inserted above the existing
# top row of the "Trading History" log, pushing the prior entry down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trading History")

# Insert a new blank row above the current top data row (row 5). Excel
# copies the header row's full column span (A:N, W:AB) and bold styling
# down onto the new row; strip that back out so only the columns actually
# used by a data row (A:J) remain, unstyled like the rest of the log.
$ws.Rows.Item(5).Insert()
$ws.Range("K5:N5").Clear()
$ws.Range("W5:AB5").Clear()
$ws.Range("A5:J5").ClearFormats()

# New (most recent) contract note entry.
$ws.Range("A5").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A5").Value = 46063
$ws.Range("B5").Value = "NSE"
$ws.Range("C5").Value = "Buy"
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 788.5
$ws.Range("F5").Value = 3970.4
$ws.Range("G5").Value = "CN#252611730667"
$ws.Range("H5").Value = 3.93
$ws.Range("I5").Value = 23.97
$ws.Range("J5").Formula = '=Index!$C$2'
